$wb = $excel.ActiveWorkbook

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1272.4565
$ws.Range("I127").Value = 2205.7
$ws.Range("K127").Value = 6617.099999999999
$ws.Range("M127").Value = -1657.099999999999

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1796.75
$ws.Range("I131").Value = 725.9231
$ws.Range("K131").Value = 2177.7693
$ws.Range("M131").Value = 2862.2307

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 978.3077
$ws.Range("I137").Value = 728.05884
$ws.Range("K137").Value = 2184.17652
$ws.Range("M137").Value = 365.82348

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2839.7827
$ws.Range("I138").Value = 2839.7827
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 8519.348100000001
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3379.348100000001
$ws.Range("N138").ClearContents()

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3661.1428
$ws.Range("I32").Value = 3269.9524
$ws.Range("J32").Value = 6008.2856
$ws.Range("K32").Value = 3269.9524
$ws.Range("L32").Value = 6008.2856
$ws.Range("M32").Value = -2982.9524
$ws.Range("N32").Value = -6582.2856

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1643.6
$ws.Range("I45").Value = 1488
$ws.Range("J45").Value = 1799.2
$ws.Range("K45").Value = 1488
$ws.Range("L45").Value = 1799.2
$ws.Range("M45").Value = -1111
$ws.Range("N45").Value = -2553.2

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2649.9143
$ws.Range("I61").Value = 1829.8572
$ws.Range("J61").Value = 5930.143
$ws.Range("K61").Value = 1829.8572
$ws.Range("L61").Value = 5930.143
$ws.Range("M61").Value = -1617.8572
$ws.Range("N61").Value = -6354.143

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1107.7106
$ws.Range("I74").Value = 916.61536
$ws.Range("J74").Value = 1521.75
$ws.Range("K74").Value = 916.61536
$ws.Range("L74").Value = 1521.75
$ws.Range("M74").Value = -42.61536000000001
$ws.Range("N74").Value = -3269.75

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1107.7106
$ws.Range("I77").Value = 916.61536
$ws.Range("J77").Value = 1521.75
$ws.Range("K77").Value = 4583.0768
$ws.Range("L77").Value = 7608.75
$ws.Range("M77").Value = -215.0767999999998
$ws.Range("N77").Value = -16344.75

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2386.182
$ws.Range("I122").Value = 2134.8
$ws.Range("K122").Value = 6404.400000000001
$ws.Range("M122").Value = -3954.400000000001

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1912.7059
$ws.Range("I132").Value = 1367.7273
$ws.Range("J132").Value = 2911.8333
$ws.Range("K132").Value = 4103.1819
$ws.Range("L132").Value = 8735.499899999999
$ws.Range("M132").Value = -1573.1819
$ws.Range("N132").Value = -13795.4999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2649.9143
$ws.Range("I136").Value = 1829.8572
$ws.Range("J136").Value = 5930.143
$ws.Range("K136").Value = 5489.571599999999
$ws.Range("L136").Value = 17790.429
$ws.Range("M136").Value = -2939.571599999999
$ws.Range("N136").Value = -22890.429

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 926.5
$ws.Range("I64").Value = 935.3333
$ws.Range("J64").Value = 900
$ws.Range("K64").Value = 935.3333
$ws.Range("L64").Value = 900
$ws.Range("M64").Value = -710.3333
$ws.Range("N64").Value = -1350

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 926.5
$ws.Range("I67").Value = 935.3333
$ws.Range("J67").Value = 900
$ws.Range("K67").Value = 935.3333
$ws.Range("L67").Value = 900
$ws.Range("M67").Value = -155.3333
$ws.Range("N67").Value = -2460

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2534.3333
$ws.Range("I105").Value = 2476.125
$ws.Range("K105").Value = 2476.125
$ws.Range("M105").Value = -729.125

# BSM row 122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 68000
$ws.Range("J122").Value = 68000
$ws.Range("L122").Value = 68000
$ws.Range("N122").Value = -77800

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6065.9395
$ws.Range("I134").Value = 6420.643
$ws.Range("J134").Value = 4079.6
$ws.Range("K134").Value = 19261.929
$ws.Range("L134").Value = 12238.8
$ws.Range("M134").Value = -16726.929
$ws.Range("N134").Value = -17308.8

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1881.0834
$ws.Range("I31").Value = 1536.3846
$ws.Range("J31").Value = 2288.4546
$ws.Range("K31").Value = 1536.3846
$ws.Range("L31").Value = 2288.4546
$ws.Range("M31").Value = -1241.3846
$ws.Range("N31").Value = -2878.4546

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1881.0834
$ws.Range("I34").Value = 1536.3846
$ws.Range("J34").Value = 2288.4546
$ws.Range("K34").Value = 1536.3846
$ws.Range("L34").Value = 2288.4546
$ws.Range("M34").Value = -1334.3846
$ws.Range("N34").Value = -2692.4546

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19812.38
$ws.Range("I131").Value = 381.14285
$ws.Range("J131").Value = 23698.629
$ws.Range("K131").Value = 1143.42855
$ws.Range("L131").Value = 71095.887
$ws.Range("M131").Value = 3896.57145
$ws.Range("N131").Value = -81175.887

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1129.2142
$ws.Range("I122").Value = 984.0833
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2952.2499
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -502.2498999999998
$ws.Range("N122").Value = -10900

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2461225.2
$ws.Range("I126").Value = 7938883
$ws.Range("J126").Value = 64749.938
$ws.Range("K126").Value = 23816649
$ws.Range("L126").Value = 194249.814
$ws.Range("M126").Value = -23814179
$ws.Range("N126").Value = -199189.814

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2164.4167
$ws.Range("I22").Value = 2082
$ws.Range("J22").Value = 2279.8
$ws.Range("K22").Value = 2082
$ws.Range("L22").Value = 2279.8
$ws.Range("M22").Value = -1787
$ws.Range("N22").Value = -2869.8

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2164.4167
$ws.Range("I27").Value = 2082
$ws.Range("J27").Value = 2279.8
$ws.Range("K27").Value = 2082
$ws.Range("L27").Value = 2279.8
$ws.Range("M27").Value = -1975
$ws.Range("N27").Value = -2493.8

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2279.0588
$ws.Range("I61").Value = 1967.4286
$ws.Range("J61").Value = 3733.3333
$ws.Range("K61").Value = 1967.4286
$ws.Range("L61").Value = 3733.3333
$ws.Range("M61").Value = -1765.4286
$ws.Range("N61").Value = -4137.3333

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1081.7333
$ws.Range("I93").Value = 786.5
$ws.Range("K93").Value = 786.5
$ws.Range("M93").Value = 461.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2279.0588
$ws.Range("I113").Value = 1967.4286
$ws.Range("J113").Value = 3733.3333
$ws.Range("K113").Value = 1967.4286
$ws.Range("L113").Value = 3733.3333
$ws.Range("M113").Value = 202.5714
$ws.Range("N113").Value = -8073.3333

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1615.0159
$ws.Range("I132").Value = 1442.5
$ws.Range("J132").Value = 1960.0476
$ws.Range("K132").Value = 4327.5
$ws.Range("L132").Value = 5880.142800000001
$ws.Range("M132").Value = -1797.5
$ws.Range("N132").Value = -10940.1428

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 70326
$ws.Range("J133").Value = 70326
$ws.Range("L133").Value = 70326
$ws.Range("N133").Value = -75386

Write-Host "Applied all edits"